# Update "想去人数" (want-to-go count) values in column F across all sheets,
# reflecting newly scraped numbers from the gh-pages data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 148
$ws1.Range("F4").Value = 602
$ws1.Range("F5").Value = 2959
$ws1.Range("F6").Value = 794
$ws1.Range("F8").Value = 601
$ws1.Range("F9").Value = 432
$ws1.Range("F12").Value = 533
$ws1.Range("F17").Value = 20
$ws1.Range("F19").Value = 2677
$ws1.Range("F23").Value = 530
$ws1.Range("F25").Value = 614
$ws1.Range("F26").Value = 13
$ws1.Range("F27").Value = 22
$ws1.Range("F31").Value = 227
$ws1.Range("F32").Value = 121
$ws1.Range("F33").Value = 901
$ws1.Range("F34").Value = 4678
$ws1.Range("F35").Value = 246
$ws1.Range("F36").Value = 30

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 10
$ws2.Range("F8").Value = 334
$ws2.Range("F9").Value = 354
$ws2.Range("F25").Value = 301
$ws2.Range("F27").Value = 161
$ws2.Range("F31").Value = 25
$ws2.Range("F36").Value = 547

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 572
$ws3.Range("F6").Value = 253
$ws3.Range("F7").Value = 258

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 572
$ws4.Range("F5").Value = 148
$ws4.Range("F6").Value = 253
$ws4.Range("F8").Value = 602
$ws4.Range("F9").Value = 2959
$ws4.Range("F10").Value = 794
$ws4.Range("F12").Value = 601
$ws4.Range("F13").Value = 432
$ws4.Range("F16").Value = 533
$ws4.Range("F17").Value = 334
$ws4.Range("F18").Value = 354
$ws4.Range("F27").Value = 2677
$ws4.Range("F32").Value = 530
$ws4.Range("F35").Value = 258
$ws4.Range("F37").Value = 614
$ws4.Range("F38").Value = 614
$ws4.Range("F42").Value = 301
$ws4.Range("F43").Value = 227
$ws4.Range("F45").Value = 901
$ws4.Range("F47").Value = 4678
$ws4.Range("F48").Value = 246
$ws4.Range("F50").Value = 547
